$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "hoodies "
$ws.Range("A4").Value = "jackets "
$ws.Range("A5").Value = "Pants"
$ws.Range("A6").Value = "Shorts "
$ws.Range("A7").Value = "Sweatshirts "
$ws.Range("A8").Value = "Tees"
$ws.Range("A9").Value = "Bras"
$ws.Range("A10").Value = "Tanks"

$ws.Range("A11").Select()
